$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("covars")

# Rename the "house weighted (standardized)" row's column_name from hv005_wi -> hiv05_wi
$ws.Range("A5").Value = "hiv05_wi"

# Select the row to be removed (mirrors selecting the whole row before deleting it in Excel)
$null = $ws.Rows.Item(9).Select()

# Drop the "hvyrmnth_dtmnth_lag" / "survey month year lagged" row (subset to de jure,
# lagged survey month/year is no longer needed) - shifts the table + everything below up by one row
$ws.Rows.Item(9).Delete()
